# Apply scheduled-runner price/profit updates to the Sheets workbook.
# Generated from the canonical OOXML diff: updates columns H-N for specific
# rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 23980.637
$ws.Range("I28").Value = 34574.24
$ws.Range("J28").Value = 5441.8335
$ws.Range("K28").Value = 34574.24
$ws.Range("L28").Value = 5441.8335
$ws.Range("M28").Value = -34089.24
$ws.Range("N28").Value = -6411.8335
$ws.Range("H38").Value = 602.4
$ws.Range("I38").Value = 8
$ws.Range("J38").Value = 857.1429000000001
$ws.Range("K38").Value = 24
$ws.Range("L38").Value = 2571.4287
$ws.Range("M38").Value = 348
$ws.Range("N38").Value = -3315.4287
$ws.Range("H113").Value = 3889.5557
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4334.3335
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4334.3335
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -10842.3335
$ws.Range("H116").Value = 3444.5
$ws.Range("I116").Value = 3709.2856
$ws.Range("K116").Value = 3709.2856
$ws.Range("M116").Value = -267.2856000000002
$ws.Range("H132").Value = 1900.3088
$ws.Range("I132").Value = 1038.75
$ws.Range("J132").Value = 3131.1072
$ws.Range("K132").Value = 3116.25
$ws.Range("L132").Value = 9393.321599999999
$ws.Range("M132").Value = -586.25
$ws.Range("N132").Value = -14453.3216
$ws.Range("H138").Value = 2004.8853
$ws.Range("I138").Value = 1121.7059
$ws.Range("J138").Value = 3117.037
$ws.Range("K138").Value = 3365.1177
$ws.Range("L138").Value = 9351.110999999999
$ws.Range("M138").Value = 1774.8823
$ws.Range("N138").Value = -19631.111
$ws.Range("H141").Value = 2816.5833
$ws.Range("I141").Value = 933.03925
$ws.Range("J141").Value = 13490
$ws.Range("K141").Value = 2799.11775
$ws.Range("L141").Value = 40470
$ws.Range("M141").Value = 2380.88225
$ws.Range("N141").Value = -50830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4422.908
$ws.Range("I32").Value = 3562.1506
$ws.Range("J32").Value = 8911.143
$ws.Range("K32").Value = 3562.1506
$ws.Range("L32").Value = 8911.143
$ws.Range("M32").Value = -3275.1506
$ws.Range("N32").Value = -9485.143
$ws.Range("H61").Value = 1071.2084
$ws.Range("I61").Value = 944.43634
$ws.Range("J61").Value = 1481.3529
$ws.Range("K61").Value = 944.43634
$ws.Range("L61").Value = 1481.3529
$ws.Range("M61").Value = -732.43634
$ws.Range("N61").Value = -1905.3529
$ws.Range("H132").Value = 5255.56
$ws.Range("I132").Value = 4079.7273
$ws.Range("K132").Value = 12239.1819
$ws.Range("M132").Value = -9709.1819
$ws.Range("H136").Value = 1071.2084
$ws.Range("I136").Value = 944.43634
$ws.Range("J136").Value = 1481.3529
$ws.Range("K136").Value = 2833.30902
$ws.Range("L136").Value = 4444.0587
$ws.Range("M136").Value = -283.3090199999997
$ws.Range("N136").Value = -9544.058700000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 56688.445
$ws.Range("I20").Value = 1178.4
$ws.Range("J20").Value = 126076
$ws.Range("K20").Value = 1178.4
$ws.Range("L20").Value = 126076
$ws.Range("M20").Value = -931.4000000000001
$ws.Range("N20").Value = -126570
$ws.Range("H94").Value = 11905595
$ws.Range("I94").Value = 20000706
$ws.Range("K94").Value = 20000706
$ws.Range("M94").Value = -20000255
$ws.Range("H105").Value = 40001492
$ws.Range("I105").Value = 62501076
$ws.Range("J105").Value = 2223.5557
$ws.Range("K105").Value = 62501076
$ws.Range("L105").Value = 2223.5557
$ws.Range("M105").Value = -62499329
$ws.Range("N105").Value = -5717.5557
$ws.Range("H134").Value = 778.05554
$ws.Range("I134").Value = 737.04346
$ws.Range("J134").Value = 1721.3334
$ws.Range("K134").Value = 2211.13038
$ws.Range("L134").Value = 5164.0002
$ws.Range("M134").Value = 323.8696199999999
$ws.Range("N134").Value = -10234.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29772.393
$ws.Range("I31").Value = 2801.1887
$ws.Range("J31").Value = 84752.16
$ws.Range("K31").Value = 2801.1887
$ws.Range("L31").Value = 84752.16
$ws.Range("M31").Value = -2506.1887
$ws.Range("N31").Value = -85342.16
$ws.Range("H34").Value = 29772.393
$ws.Range("I34").Value = 2801.1887
$ws.Range("J34").Value = 84752.16
$ws.Range("K34").Value = 2801.1887
$ws.Range("L34").Value = 84752.16
$ws.Range("M34").Value = -2599.1887
$ws.Range("N34").Value = -85156.16
$ws.Range("H99").Value = 3105.5
$ws.Range("I99").Value = 1704
$ws.Range("J99").Value = 4507
$ws.Range("K99").Value = 1704
$ws.Range("L99").Value = 4507
$ws.Range("M99").Value = -206
$ws.Range("N99").Value = -7503
$ws.Range("H126").Value = 3105.5
$ws.Range("I126").Value = 1704
$ws.Range("J126").Value = 4507
$ws.Range("K126").Value = 5112
$ws.Range("L126").Value = 13521
$ws.Range("M126").Value = -2642
$ws.Range("N126").Value = -18461
$ws.Range("H132").Value = 16397027
$ws.Range("I132").Value = 23259878
$ws.Range("J132").Value = 2439.6667
$ws.Range("K132").Value = 69779634
$ws.Range("L132").Value = 7319.000100000001
$ws.Range("M132").Value = -69777104
$ws.Range("N132").Value = -12379.0001
$ws.Range("H134").Value = 3631.2888
$ws.Range("I134").Value = 4157
$ws.Range("J134").Value = 1528.4445
$ws.Range("K134").Value = 12471
$ws.Range("L134").Value = 4585.333500000001
$ws.Range("M134").Value = -9936
$ws.Range("N134").Value = -9655.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 707.11365
$ws.Range("J131").Value = 992.2222
$ws.Range("L131").Value = 2976.6666
$ws.Range("N131").Value = -13056.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3630.377
$ws.Range("I132").Value = 4071.182
$ws.Range("J132").Value = 2489.4707
$ws.Range("K132").Value = 12213.546
$ws.Range("L132").Value = 7468.4121
$ws.Range("M132").Value = -9683.545999999998
$ws.Range("N132").Value = -12528.4121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3757.0645
$ws.Range("I132").Value = 4670.162
$ws.Range("J132").Value = 2405.68
$ws.Range("K132").Value = 14010.486
$ws.Range("L132").Value = 7217.039999999999
$ws.Range("M132").Value = -11480.486
$ws.Range("N132").Value = -12277.04

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2134.353
$ws.Range("I132").Value = 2980.0217
$ws.Range("J132").Value = 1136.8975
$ws.Range("K132").Value = 8940.0651
$ws.Range("L132").Value = 3410.6925
$ws.Range("M132").Value = -6410.0651
$ws.Range("N132").Value = -8470.692500000001
$ws.Range("H136").Value = 1699.0793
$ws.Range("I136").Value = 1911.8572
$ws.Range("J136").Value = 1528.8572
$ws.Range("K136").Value = 5735.571599999999
$ws.Range("L136").Value = 4586.571599999999
$ws.Range("M136").Value = -3185.571599999999
$ws.Range("N136").Value = -9686.571599999999
